$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1163-1164, pushing the existing data (old rows
# 1163-1259) down to 1165-1261.
$ws.Rows("1163:1164").Insert()

# Populate the two newly inserted rows with the new weekly records.
$ws.Range("A1163").Value = 6
$ws.Range("B1163").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1163").Value = "Metropolitana"
$ws.Range("D1163").Value = 45106
$ws.Range("E1163").Value = 13
$ws.Range("F1163").Value = 100112009
$ws.Range("G1163").Value = "Acelga"
$ws.Range("H1163").Value = "Sin especificar"
$ws.Range("I1163").Value = "Primera"
$ws.Range("J1163").Value = 250
$ws.Range("K1163").Value = 12000
$ws.Range("L1163").Value = 12000
$ws.Range("M1163").Value = 12000
$ws.Range("N1163").Value = "`$/docena de atados"
$ws.Range("O1163").Value = "Región Metropolitana"
$ws.Range("P1163").Value = 4000
$ws.Range("Q1163").Value = 3
$ws.Range("R1163").Value = "Hortaliza"

$ws.Range("A1164").Value = 6
$ws.Range("B1164").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1164").Value = "Metropolitana"
$ws.Range("D1164").Value = 45106
$ws.Range("E1164").Value = 13
$ws.Range("F1164").Value = 100112009
$ws.Range("G1164").Value = "Acelga"
$ws.Range("H1164").Value = "Sin especificar"
$ws.Range("I1164").Value = "Segunda"
$ws.Range("J1164").Value = 170
$ws.Range("K1164").Value = 9000
$ws.Range("L1164").Value = 9000
$ws.Range("M1164").Value = 9000
$ws.Range("N1164").Value = "`$/docena de atados"
$ws.Range("O1164").Value = "Región Metropolitana"
$ws.Range("P1164").Value = 3000
$ws.Range("Q1164").Value = 3
$ws.Range("R1164").Value = "Hortaliza"
